# Test data document name changed:
#   FLD_CreateNewDocument  ->  Documents_CreateNewDocument
# This renames the referenced test-data workbook from
# "FLD_CreateNewDocument.xlsx" to "Documents_CreateNewDocument.xlsx" and
# updates the two lookup-name cells (sheet1!A21 and sheet2!C21) that held
# the old "FLD_CreateNewDocument" label to the new
# "Documents_CreateNewDocument" label. It also moves the active
# sheet/selection from DataFetchXL (sheet2) to DataFetchFlag (sheet1), and
# updates each sheet's remembered selection to match.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("DataFetchFlag")
$ws2 = $wb.Worksheets.Item("DataFetchXL")

# --- Data edits -------------------------------------------------------
# (Order matters for how the shared-string table comes out: update the
# path text first, while its cell is still the only one pointing at the
# old "FLD_CreateNewDocument.xlsx" string, so it's rewritten in place;
# then rename the two label cells so the new "Documents_CreateNewDocument"
# label is added once and shared.)

# sheet2 (DataFetchXL) B21: update the remembered file path text to match
# the renamed test-data document (trailing space preserved, as in source).
$ws2.Range("B21").Value = "\\src\com\proj\suiteDOCS\testdata\Documents_CreateNewDocument.xlsx "

# sheet1 (DataFetchFlag) A21: FLD_CreateNewDocument -> Documents_CreateNewDocument
$ws1.Range("A21").Value = "Documents_CreateNewDocument"

# sheet2 (DataFetchXL) C21: FLD_CreateNewDocument -> Documents_CreateNewDocument
$ws2.Range("C21").Value = "Documents_CreateNewDocument"

# --- View-state edits ---------------------------------------------------

# Sheet2's selection moves to C21.
$ws2.Activate()
$ws2.Range("C21").Select()

# Sheet1 becomes the active tab/sheet; selection moves to A21. (Selecting
# on sheet2 above switches the active sheet, so sheet1 is re-activated
# last so it ends up as the active/selected tab, matching the target.)
$ws1.Activate()
$ws1.Range("A21").Select()
